# Generate Report for Handoff
# The f6741cfc-a833-45ab-b1a9-110e314812dd.md file (row 3 in each sheet) is
# moving from "In Translation" to "Ready for handoff". Update the Overview
# sheet's rollup row plus the per-locale (zh-cn / de-de) detail rows with
# their new status, priority and handoff timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-21 08:13:41"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-21 08:13:37"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-21 08:13:41"
